# Insert a new log row (2026/01/10, 土, 6, 25) above the current row 591
# ("2026/12/29" row), pushing that row and everything below it down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(591).Insert()

# Copy the formatting of the date cell directly above (which already holds
# the plain-text date "2026/01/10") onto the new A591 so the freshly
# inserted row doesn't pick up a stray number/date style.
$ws.Cells.Item(590, 1).Copy()
$ws.Cells.Item(591, 1).PasteSpecial()

$ws.Cells.Item(591, 2).Value = "土"
$ws.Cells.Item(591, 3).Value = 6
$ws.Cells.Item(591, 4).Value = 25
